$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates scraped from the "Updated symbol list" refresh: each
# row gets a fresh Price (D), Volume(1h) (E) and Hora (G) reading. A handful
# of rows (placeholder "--" coins, and CoinbaseStockToken) only had their
# Hora bumped, so D/E are $null there and left untouched below.
$rows = @(
    @{ Row=2; D="320.64"; E="-3.35%"; G="13" },
    @{ Row=3; D="42.87"; E="-5.67%"; G="13" },
    @{ Row=4; D="5.204"; E="-5.47%"; G="13" },
    @{ Row=5; D="0.08158"; E="-3.63%"; G="13" },
    @{ Row=6; D="4.325"; E="-2.35%"; G="13" },
    @{ Row=7; D="1.802"; E="-14.14%"; G="13" },
    @{ Row=8; D="0.9417"; E="-4.88%"; G="13" },
    @{ Row=9; D="0.1112"; E="-6.92%"; G="13" },
    @{ Row=10; D="0.1859"; E="-3.64%"; G="13" },
    @{ Row=11; D="0.09367"; E="-4.44%"; G="13" },
    @{ Row=12; D="0.04630"; E="-1.06%"; G="13" },
    @{ Row=13; D="7.430"; E="-22.34%"; G="13" },
    @{ Row=14; D="0.1058"; E="-0.26%"; G="13" },
    @{ Row=15; D="0.001292"; E="-0.81%"; G="13" },
    @{ Row=16; D="0.005949"; E="0.46%"; G="13" },
    @{ Row=17; D="3.361"; E="-0.83%"; G="13" },
    @{ Row=18; D="2.548"; E="-0.24%"; G="13" },
    @{ Row=19; D="0.3349"; E="0.32%"; G="13" },
    @{ Row=20; D="0.1380"; E="1.80%"; G="13" },
    @{ Row=21; D="0.2550"; E="0.07%"; G="13" },
    @{ Row=22; D="0.04144"; E="-0.05%"; G="13" },
    @{ Row=23; D="0.001243"; E="-4.42%"; G="13" },
    @{ Row=24; D="0.004277"; E="-6.12%"; G="13" },
    @{ Row=25; D="0.0001201"; E="-7.93%"; G="13" },
    @{ Row=26; D="0.0002980"; E="-20.44%"; G="13" },
    @{ Row=27; D=$null; E=$null; G="13" },
    @{ Row=28; D=$null; E=$null; G="13" },
    @{ Row=29; D=$null; E=$null; G="13" },
    @{ Row=30; D=$null; E=$null; G="13" },
    @{ Row=31; D=$null; E=$null; G="13" },
    @{ Row=32; D=$null; E=$null; G="13" },
    @{ Row=33; D=$null; E=$null; G="13" },
    @{ Row=34; D=$null; E=$null; G="13" },
    @{ Row=35; D=$null; E=$null; G="13" },
    @{ Row=36; D=$null; E=$null; G="13" },
    @{ Row=37; D=$null; E=$null; G="13" },
    @{ Row=38; D="0.02689"; E="-0.23%"; G="13" },
    @{ Row=39; D="0.05548"; E="-4.13%"; G="13" },
    @{ Row=40; D="0.008096"; E="4.42%"; G="13" },
    @{ Row=41; D="0.1401"; E="-2.23%"; G="13" },
    @{ Row=42; D="0.006558"; E="-12.79%"; G="13" },
    @{ Row=43; D="0.002086"; E="-0.97%"; G="13" },
    @{ Row=44; D="0.007599"; E="-15.15%"; G="13" },
    @{ Row=45; D="0.3186"; E="-10.11%"; G="13" },
    @{ Row=46; D="0.00006928"; E="-2.87%"; G="13" },
    @{ Row=47; D="0.00000000750"; E="-0.25%"; G="13" },
    @{ Row=48; D="0.003335"; E="10.93%"; G="13" },
    @{ Row=49; D=$null; E="-0.27%"; G="13" },
    @{ Row=50; D="0.00002101"; E="-0.25%"; G="13" },
    @{ Row=51; D="0.0002001"; E="-0.25%"; G="13" }
)

foreach ($r in $rows) {
    # Values like "320.64", "-3.35%" and "13" look numeric/percent to Excel,
    # which would silently coerce them to numbers on assignment. Forcing
    # Text format on each target cell first keeps them as literal strings.
    if ($null -ne $r.D) {
        $ws.Range("D" + $r.Row).NumberFormat = "@"
        $ws.Range("D" + $r.Row).Value = $r.D
    }
    if ($null -ne $r.E) {
        $ws.Range("E" + $r.Row).NumberFormat = "@"
        $ws.Range("E" + $r.Row).Value = $r.E
    }
    $ws.Range("G" + $r.Row).NumberFormat = "@"
    $ws.Range("G" + $r.Row).Value = $r.G
}
